$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply new averaged-intensity data for rows 10-19 (HKL indices 8-17).
# Existing rows 10-16 are being overwritten in place (same row numbers,
# the underlying scheme each row represents has shifted because three new
# "Spiral" rotation schemes were measured and inserted), and rows 17-19 are
# brand new rows appended at the bottom of the table.

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.8620159336131725
$ws.Range("D10").Value = 1.1851820878893
$ws.Range("E10").Value = 0.9528535474453399
$ws.Range("F10").Value = 1.036449921668027
$ws.Range("G10").Value = 0.8620159336131725
$ws.Range("H10").Value = 1.1851820878893
$ws.Range("I10").Value = 0.9329395222185278
$ws.Range("J10").Value = 1.033450248940469
$ws.Range("K10").Value = 0.9489456796225956
$ws.Range("L10").Value = 1.11605938060394
$ws.Range("M10").Value = 0.8620159336131725
$ws.Range("N10").Value = 1.06901781766732
$ws.Range("O10").Value = 1.00912537265396
$ws.Range("P10").Value = 1.008487040250172

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.7955389212182482
$ws.Range("D11").Value = 1.268302539030818
$ws.Range("E11").Value = 0.9440334041316302
$ws.Range("F11").Value = 1.053479590065822
$ws.Range("G11").Value = 0.7955389212182482
$ws.Range("H11").Value = 1.268302539030818
$ws.Range("I11").Value = 0.9017252747761558
$ws.Range("J11").Value = 1.059363709310857
$ws.Range("K11").Value = 0.9210012378633745
$ws.Range("L11").Value = 1.170050489518952
$ws.Range("M11").Value = 0.7955389212182482
$ws.Range("N11").Value = 1.106167971581224
$ws.Range("O11").Value = 1.01533861361163
$ws.Range("P11").Value = 1.014186895739482

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.7963266716619661
$ws.Range("D12").Value = 1.266842653262206
$ws.Range("E12").Value = 0.9444888236289528
$ws.Range("F12").Value = 1.053134455827418
$ws.Range("G12").Value = 0.7963266716619661
$ws.Range("H12").Value = 1.266842653262206
$ws.Range("I12").Value = 0.9022533326309756
$ws.Range("J12").Value = 1.059198500077432
$ws.Range("K12").Value = 0.9212553587219284
$ws.Range("L12").Value = 1.169164959788629
$ws.Range("M12").Value = 0.7963266716619661
$ws.Range("N12").Value = 1.105665738445579
$ws.Range("O12").Value = 1.015198151095136
$ws.Range("P12").Value = 1.014083094449939

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.7957305428138802
$ws.Range("D13").Value = 1.268034431053711
$ws.Range("E13").Value = 0.9440929735483513
$ws.Range("F13").Value = 1.053371803780082
$ws.Range("G13").Value = 0.7957305428138802
$ws.Range("H13").Value = 1.268034431053711
$ws.Range("I13").Value = 0.9018274802579159
$ws.Range("J13").Value = 1.059351485497782
$ws.Range("K13").Value = 0.9210400306190084
$ws.Range("L13").Value = 1.169911622423637
$ws.Range("M13").Value = 0.7957305428138802
$ws.Range("N13").Value = 1.106063702301031
$ws.Range("O13").Value = 1.015307437799006
$ws.Range("P13").Value = 1.014170046249296

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.6694800000000011
$ws.Range("D14").Value = 1.510803999999997
$ws.Range("E14").Value = 0.8643160000000008
$ws.Range("F14").Value = 1.112967999999999
$ws.Range("G14").Value = 0.6694800000000011
$ws.Range("H14").Value = 1.510803999999997
$ws.Range("I14").Value = 0.8170079999999995
$ws.Range("J14").Value = 1.084396
$ws.Range("K14").Value = 0.8780919999999999
$ws.Range("L14").Value = 1.301987999999999
$ws.Range("M14").Value = 0.6694800000000011
$ws.Range("N14").Value = 1.187559999999999
$ws.Range("O14").Value = 1.039391999999999
$ws.Range("P14").Value = 1.0298815

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 1.830412500000001
$ws.Range("E15").Value = 0.76
$ws.Range("F15").Value = 1.188487499999999
$ws.Range("G15").Value = 0.5
$ws.Range("H15").Value = 1.830412500000001
$ws.Range("I15").Value = 0.7
$ws.Range("J15").Value = 1.12
$ws.Range("K15").Value = 0.8187625000000001
$ws.Range("L15").Value = 1.49
$ws.Range("M15").Value = 0.5
$ws.Range("N15").Value = 1.29520625
$ws.Range("O15").Value = 1.069725
$ws.Range("P15").Value = 1.0509578125

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.7089890035712006
$ws.Range("D16").Value = 1.481975428608
$ws.Range("E16").Value = 0.8566873790464007
$ws.Range("F16").Value = 1.105847914700799
$ws.Range("G16").Value = 0.7089890035712006
$ws.Range("H16").Value = 1.481975428608
$ws.Range("I16").Value = 0.8239816525823993
$ws.Range("J16").Value = 1.069126881484794
$ws.Range("K16").Value = 0.8915117553664049
$ws.Range("L16").Value = 1.283221288755201
$ws.Range("M16").Value = 0.7089890035712006
$ws.Range("N16").Value = 1.1693314038272
$ws.Range("O16").Value = 1.0383749314816
$ws.Range("P16").Value = 1.0276676630144

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9936995693683933
$ws.Range("D17").Value = 0.9951408807131217
$ws.Range("E17").Value = 0.9966757790444358
$ws.Range("F17").Value = 0.9965978855928875
$ws.Range("G17").Value = 0.9936995693683933
$ws.Range("H17").Value = 0.9951408807131217
$ws.Range("I17").Value = 0.9964990863024565
$ws.Range("J17").Value = 0.9946632626008348
$ws.Range("K17").Value = 0.9951164975163977
$ws.Range("L17").Value = 0.997713137552179
$ws.Range("M17").Value = 0.9936995693683933
$ws.Range("N17").Value = 0.9959083298787788
$ws.Range("O17").Value = 0.9955285286797095
$ws.Range("P17").Value = 0.9957632623363383

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.013051849414023
$ws.Range("D18").Value = 0.9683805076746137
$ws.Range("E18").Value = 1.001460038053012
$ws.Range("F18").Value = 0.9913880389249421
$ws.Range("G18").Value = 1.013051849414023
$ws.Range("H18").Value = 0.9683805076746137
$ws.Range("I18").Value = 1.007297022730692
$ws.Range("J18").Value = 0.9879392069887653
$ws.Range("K18").Value = 1.002788084881272
$ws.Range("L18").Value = 0.9797339078626477
$ws.Range("M18").Value = 1.013051849414023
$ws.Range("N18").Value = 0.9849202728638127
$ws.Range("O18").Value = 0.9935701085166477
$ws.Range("P18").Value = 0.994004832066246

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.054633662475778
$ws.Range("D19").Value = 0.9096925873334168
$ws.Range("E19").Value = 1.014809039204192
$ws.Range("F19").Value = 0.976056732362078
$ws.Range("G19").Value = 1.054633662475778
$ws.Range("H19").Value = 0.9096925873334168
$ws.Range("I19").Value = 1.029167065823824
$ws.Range("J19").Value = 0.9787548631096405
$ws.Range("K19").Value = 1.015627061978553
$ws.Range("L19").Value = 0.9380930546685059
$ws.Range("M19").Value = 1.054633662475778
$ws.Range("N19").Value = 0.9622508132688044
$ws.Range("O19").Value = 0.9887980053438663
$ws.Range("P19").Value = 0.9896042583694986

# New rows 17-19 need the same bold/bordered/centered style as the rest of
# column A (the HKL-index column) -- copy formatting from an existing styled cell.
$ws.Range("A10").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0
